$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap match data (cols F:V) between rows 15 and 16
$ws.Range("F15").Value = "Galatasaray"
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = "Trabzonspor"
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 1.44
$ws.Range("K15").Value = "13/08/2023 01:12"
$ws.Range("L15").Value = 1.76
$ws.Range("M15").Value = "19/08/2023 20:44"
$ws.Range("N15").Value = 5.28
$ws.Range("O15").Value = "13/08/2023 01:12"
$ws.Range("P15").Value = 4.14
$ws.Range("Q15").Value = "19/08/2023 20:44"
$ws.Range("R15").Value = 6.67
$ws.Range("S15").Value = "13/08/2023 01:12"
$ws.Range("T15").Value = 4.51
$ws.Range("U15").Value = "19/08/2023 20:44"
$ws.Range("V15").Value = "https://www.betexplorer.com/football/turkey/super-lig/galatasaray-trabzonspor/GdlCZEXT/"
$ws.Range("F16").Value = "Hatayspor"
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = "Kasimpasa"
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 2.54
$ws.Range("K16").Value = "13/08/2023 09:12"
$ws.Range("L16").Value = 2.63
$ws.Range("M16").Value = "19/08/2023 20:43"
$ws.Range("N16").Value = 3.59
$ws.Range("O16").Value = "13/08/2023 09:12"
$ws.Range("P16").Value = 3.43
$ws.Range("Q16").Value = "19/08/2023 17:05"
$ws.Range("R16").Value = 2.79
$ws.Range("S16").Value = "13/08/2023 09:12"
$ws.Range("T16").Value = 2.85
$ws.Range("U16").Value = "19/08/2023 20:09"
$ws.Range("V16").Value = "https://www.betexplorer.com/football/turkey/super-lig/hatayspor-kasimpasa/rVz5yy3H/"

# Swap match data (cols F:V) between rows 18 and 19
$ws.Range("F18").Value = "Besiktas"
$ws.Range("G18").Value = 1
$ws.Range("H18").Value = "Pendikspor"
$ws.Range("I18").Value = 1
$ws.Range("J18").Value = 1.29
$ws.Range("K18").Value = "15/08/2023 13:42"
$ws.Range("L18").Value = 1.28
$ws.Range("M18").Value = "20/08/2023 20:12"
$ws.Range("N18").Value = 6.05
$ws.Range("O18").Value = "15/08/2023 13:42"
$ws.Range("P18").Value = 6.31
$ws.Range("Q18").Value = "20/08/2023 20:12"
$ws.Range("R18").Value = 8.81
$ws.Range("S18").Value = "15/08/2023 13:42"
$ws.Range("T18").Value = 10.07
$ws.Range("U18").Value = "20/08/2023 20:12"
$ws.Range("V18").Value = "https://www.betexplorer.com/football/turkey/super-lig/besiktas-pendikspor/6RGYjG24/"
$ws.Range("F19").Value = "Gaziantep"
$ws.Range("G19").Value = 1
$ws.Range("H19").Value = "Sivasspor"
$ws.Range("I19").Value = 3
$ws.Range("J19").Value = 2.45
$ws.Range("K19").Value = "15/08/2023 13:42"
$ws.Range("L19").Value = 2.84
$ws.Range("M19").Value = "20/08/2023 20:43"
$ws.Range("N19").Value = 3.59
$ws.Range("O19").Value = "15/08/2023 13:42"
$ws.Range("P19").Value = 3.42
$ws.Range("Q19").Value = "20/08/2023 20:43"
$ws.Range("R19").Value = 2.9
$ws.Range("S19").Value = "15/08/2023 13:42"
$ws.Range("T19").Value = 2.61
$ws.Range("U19").Value = "20/08/2023 20:43"
$ws.Range("V19").Value = "https://www.betexplorer.com/football/turkey/super-lig/gaziantep-sivasspor/KAEMgIIo/"

# Swap match data (cols F:V) between rows 26 and 27
$ws.Range("F26").Value = "Kayserispor"
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = "Samsunspor"
$ws.Range("I26").Value = 1
$ws.Range("J26").Value = 2.32
$ws.Range("K26").Value = "22/08/2023 07:12"
$ws.Range("L26").Value = 2.68
$ws.Range("M26").Value = "27/08/2023 17:46"
$ws.Range("N26").Value = 3.62
$ws.Range("O26").Value = "22/08/2023 07:12"
$ws.Range("P26").Value = 3.35
$ws.Range("Q26").Value = "27/08/2023 18:11"
$ws.Range("R26").Value = 3.09
$ws.Range("S26").Value = "22/08/2023 07:12"
$ws.Range("T26").Value = 2.81
$ws.Range("U26").Value = "27/08/2023 18:11"
$ws.Range("V26").Value = "https://www.betexplorer.com/football/turkey/super-lig/kayserispor-samsunspor/fkYAXR1q/"
$ws.Range("F27").Value = "Karagumruk"
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = "Ankaragucu"
$ws.Range("I27").Value = 1
$ws.Range("J27").Value = 2.51
$ws.Range("K27").Value = "21/08/2023 20:12"
$ws.Range("L27").Value = 2.48
$ws.Range("M27").Value = "27/08/2023 18:13"
$ws.Range("N27").Value = 3.65
$ws.Range("O27").Value = "21/08/2023 20:12"
$ws.Range("P27").Value = 3.38
$ws.Range("Q27").Value = "27/08/2023 18:12"
$ws.Range("R27").Value = 2.78
$ws.Range("S27").Value = "21/08/2023 20:12"
$ws.Range("T27").Value = 3.05
$ws.Range("U27").Value = "27/08/2023 18:13"
$ws.Range("V27").Value = "https://www.betexplorer.com/football/turkey/super-lig/f-karagumruk-ankaragucu/foI9hTw9/"

# Swap match data (cols F:V) between rows 51 and 52
$ws.Range("F51").Value = "Basaksehir"
$ws.Range("G51").Value = 1
$ws.Range("H51").Value = "Galatasaray"
$ws.Range("I51").Value = 2
$ws.Range("J51").Value = 5.38
$ws.Range("K51").Value = "17/09/2023 15:12"
$ws.Range("L51").Value = 4.24
$ws.Range("M51").Value = "23/09/2023 18:57"
$ws.Range("N51").Value = 4.5
$ws.Range("O51").Value = "17/09/2023 15:12"
$ws.Range("P51").Value = 3.75
$ws.Range("Q51").Value = "23/09/2023 18:53"
$ws.Range("R51").Value = 1.6
$ws.Range("S51").Value = "17/09/2023 15:12"
$ws.Range("T51").Value = 1.9
$ws.Range("U51").Value = "23/09/2023 18:57"
$ws.Range("V51").Value = "https://www.betexplorer.com/football/turkey/super-lig/basaksehir-galatasaray/KjY6EVHP/"
$ws.Range("F52").Value = "Antalyaspor"
$ws.Range("G52").Value = 2
$ws.Range("H52").Value = "Samsunspor"
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 1.97
$ws.Range("K52").Value = "19/09/2023 14:42"
$ws.Range("L52").Value = 2.18
$ws.Range("M52").Value = "23/09/2023 18:52"
$ws.Range("N52").Value = 3.8
$ws.Range("O52").Value = "19/09/2023 14:42"
$ws.Range("P52").Value = 3.55
$ws.Range("Q52").Value = "23/09/2023 18:52"
$ws.Range("R52").Value = 3.82
$ws.Range("S52").Value = "19/09/2023 14:42"
$ws.Range("T52").Value = 3.47
$ws.Range("U52").Value = "23/09/2023 18:52"
$ws.Range("V52").Value = "https://www.betexplorer.com/football/turkey/super-lig/antalyaspor-samsunspor/xWKY9kQt/"

# Swap match data (cols F:V) between rows 55 and 56
$ws.Range("F55").Value = "Alanyaspor"
$ws.Range("G55").Value = 0
$ws.Range("H55").Value = "Fenerbahce"
$ws.Range("I55").Value = 1
$ws.Range("J55").Value = 5.29
$ws.Range("K55").Value = "17/09/2023 18:12"
$ws.Range("L55").Value = 5.65
$ws.Range("M55").Value = "24/09/2023 18:59"
$ws.Range("N55").Value = 4.73
$ws.Range("O55").Value = "17/09/2023 18:12"
$ws.Range("P55").Value = 4.49
$ws.Range("Q55").Value = "24/09/2023 18:59"
$ws.Range("R55").Value = 1.57
$ws.Range("S55").Value = "17/09/2023 18:12"
$ws.Range("T55").Value = 1.58
$ws.Range("U55").Value = "24/09/2023 18:59"
$ws.Range("V55").Value = "https://www.betexplorer.com/football/turkey/super-lig/alanyaspor-fenerbahce/2F3M6JbJ/"
$ws.Range("F56").Value = "Kasimpasa"
$ws.Range("G56").Value = 2
$ws.Range("H56").Value = "Adana Demirspor"
$ws.Range("I56").Value = 1
$ws.Range("J56").Value = 2.79
$ws.Range("K56").Value = "17/09/2023 18:12"
$ws.Range("L56").Value = 3.51
$ws.Range("M56").Value = "24/09/2023 18:59"
$ws.Range("N56").Value = 3.85
$ws.Range("O56").Value = "17/09/2023 18:12"
$ws.Range("P56").Value = 4.14
$ws.Range("Q56").Value = "24/09/2023 18:59"
$ws.Range("R56").Value = 2.37
$ws.Range("S56").Value = "17/09/2023 18:12"
$ws.Range("T56").Value = 1.99
$ws.Range("U56").Value = "24/09/2023 18:58"
$ws.Range("V56").Value = "https://www.betexplorer.com/football/turkey/super-lig/kasimpasa-adanademirspor/OO4I7aqD/"

# Append new match row 108 (Indice 107)
$ws.Range("A107").Copy() | Out-Null
$ws.Range("A108").PasteSpecial(-4122) | Out-Null
$ws.Range("E107").Copy() | Out-Null
$ws.Range("E108").PasteSpecial(-4122) | Out-Null

$ws.Range("A108").Value = 107
$ws.Range("B108").Value = "turkey"
$ws.Range("C108").Value = "super-lig"
$ws.Range("D108").Value = "2023-2024"
$ws.Range("E108").Value = 45235.70833333334
$ws.Range("F108").Value = "Antalyaspor"
$ws.Range("G108").Value = 3
$ws.Range("H108").Value = "Besiktas"
$ws.Range("I108").Value = 2
$ws.Range("J108").Value = 3.34
$ws.Range("K108").Value = "30/10/2023 18:12"
$ws.Range("L108").Value = 2.55
$ws.Range("M108").Value = "05/11/2023 16:59"
$ws.Range("N108").Value = 3.63
$ws.Range("O108").Value = "30/10/2023 18:12"
$ws.Range("P108").Value = 3.52
$ws.Range("Q108").Value = "05/11/2023 16:51"
$ws.Range("R108").Value = 2.15
$ws.Range("S108").Value = "30/10/2023 18:12"
$ws.Range("T108").Value = 2.85
$ws.Range("U108").Value = "05/11/2023 16:59"
$ws.Range("V108").Value = "https://www.betexplorer.com/football/turkey/super-lig/antalyaspor-besiktas/vV7GsRss/"
